$d = $word.ActiveDocument
$d.Content.Find.Execute("Molex Connected Mobility", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Molex LLC", 2)
